$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows keyed by row number
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 10132
$wsExpo.Range("F10").Value = 6590
$wsExpo.Range("F13").Value = 188
$wsExpo.Range("F15").Value = 3116
$wsExpo.Range("F23").Value = 1556

# Sheet "全部类型" (all types) - same events, rows shifted by +1
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 10132
$wsAll.Range("F11").Value = 6590
$wsAll.Range("F14").Value = 188
$wsAll.Range("F16").Value = 3116
$wsAll.Range("F24").Value = 1556
